# Exporting metrics for each file in the analysed project.
#
# This adds a new "Metrics" worksheet (mirroring the layout already used by
# the "All"/"Unconfirmed" sheets: a single "Colonne1" table fed by the
# report generator), makes it the last / active sheet, and refreshes the
# existing "synthesis" pivot table so it stays in sync with its source data.

$wb = $excel.ActiveWorkbook

# Template sheet to mimic (same single-column "Colonne1" table layout).
$template = $wb.Worksheets.Item("All")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Create the new worksheet as the last tab in the workbook.
$metrics = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$metrics.Name = "Metrics"

# Match column widths used on the sibling "All" sheet.
$metrics.Columns.Item(2).ColumnWidth = $template.Columns.Item(2).ColumnWidth
$metrics.Columns.Item(3).ColumnWidth = $template.Columns.Item(3).ColumnWidth
$metrics.Columns.Item(10).ColumnWidth = $template.Columns.Item(10).ColumnWidth
$metrics.Columns.Item(13).ColumnWidth = $template.Columns.Item(13).ColumnWidth
$metrics.Columns.Item(15).ColumnWidth = $template.Columns.Item(15).ColumnWidth
$metrics.Columns.Item(16).ColumnWidth = $template.Columns.Item(16).ColumnWidth

# Header cell/row, styled like the other "Colonne1" tables.
$metrics.Range("A1").Value = "Colonne1"
$metrics.Rows.Item(1).RowHeight = 30
$metrics.Range("A1").HorizontalAlignment = -4108
$metrics.Range("A1").VerticalAlignment = -4108

# Turn the header + first body row into a table, like "all"/"unconfirmed".
$metricsTable = $metrics.ListObjects.Add(1, $metrics.Range("A1:A2"), [System.Reflection.Missing]::Value, 1)
$metricsTable.Name = "metrics"
$metricsTable.TableStyle = "TableStyleLight16"

# Page setup consistent with the other exported-table sheets.
$metrics.PageSetup.PaperSize = 9
$metrics.PageSetup.Orientation = 1

# New sheet becomes the active / selected tab.
$metrics.Activate()

# Keep the pivot table ("synthesis" on "TCD") in sync with its source table.
$tcd = $wb.Worksheets.Item("TCD")
foreach ($pt in $tcd.PivotTables()) {
    $pt.RefreshTable()
}

Write-Output "Added Metrics sheet with metrics table"
